$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.122.08"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.175.79"
$ws.Range("E3").Value = "  +3.51%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.94"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.48"
$ws.Range("E6").Value = "  +4.60%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.175.99"
$ws.Range("E8").Value = "  +3.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.506"
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000276"
$ws.Range("E13").Value = "  +18.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.07"
$ws.Range("E14").Value = "  +6.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.696.55"
$ws.Range("E15").Value = "  +3.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.200.71"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.179.04"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  +6.33%  "
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.07"
$ws.Range("E20").Value = "  +7.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.92"
$ws.Range("E21").Value = "  +6.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("E22").Value = "  +7.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.39"
$ws.Range("E23").Value = "  +7.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.84"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.28"
$ws.Range("E25").Value = "  +3.58%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.09"
$ws.Range("E27").Value = "  +11.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.91"
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  +7.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.14"
$ws.Range("E30").Value = "  +6.57%  "
$ws.Range("E31").Value = "  +13.18%  "
$ws.Range("E32").Value = "  +6.61%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +9.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.70"
$ws.Range("E35").Value = "  +6.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.77"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0900"
$ws.Range("E37").Value = "  +10.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "480.30"
$ws.Range("E38").Value = "  +7.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").Value = "  +8.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0422"
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.115.12"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.65"
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("E43").Value = "  +3.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.288"
$ws.Range("E44").Value = "  +9.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  +12.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.44"
$ws.Range("E46").Value = "  +5.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0599"
$ws.Range("E47").Value = "  +14.84%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("E50").Value = "  +9.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.62"
$ws.Range("E51").Value = "  +2.42%  "
